$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.719595074653625
$ws.Range("B1").Value = 1.744430541992188
$ws.Range("C1").Value = 4.880401134490967
$ws.Range("D1").Value = 1.208985686302185
$ws.Range("E1").Value = 0.637040376663208
